# Update the poster date from "Dec 4th-9th, 2022" to "Dec 5th-9th, 2022".
#
# The "Duration" line lives in TextBox 6 on slide 1 and is split into
# several runs:
#   [Duration][: Dec ][4][th][-9][th][, 2022]
# The edit merges the ": Dec " run with the "4" run into a single run
# with text ": Dec 5" (keeping the formatting of the ": Dec " run),
# and removes the now-empty "4" run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTextFrame -and $candidate.TextFrame.TextRange.Text.StartsWith("Duration:")) {
        $shape = $candidate
        break
    }
}

$tr = $shape.TextFrame.TextRange

# Characters 9-15 of the text cover ": Dec " (6 chars) immediately followed
# by "4" (1 char) -- i.e. the two runs being collapsed into one.
$target = $tr.Characters(9, 7)
$target.Text = ": Dec 5"
